$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "ECO Actual" (B) and "ECO Balance" (D) columns hold numeric-looking
# values stored as text. Update them to the new reported figures while
# preserving their text formatting (e.g. "5.00" rather than 5).

$updates = @{
    "B2" = "5.00"
    "D2" = "5.00"
    "B3" = "9.00"
    "D3" = "9.00"
    "B4" = "6.00"
    "D4" = "6.00"
    "B5" = "8.00"
    "D5" = "8.00"
    "B6" = "6.00"
    "D6" = "6.00"
    "B7" = "34.00"
    "D7" = "34.00"
}

foreach ($addr in $updates.Keys) {
    # Leading apostrophe forces Excel to keep the numeric-looking text
    # (e.g. "5.00") as a literal text value instead of coercing it to
    # the number 5, matching how the original sheet stores these figures.
    $ws.Range($addr).Value = "'" + $updates[$addr]
}
